$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"9.9457514024314158"
$ws.Range("A3").Value = [double]"-16.22779182725651"
$ws.Range("B3").Value = [double]"12.872114284246949"
$ws.Range("A4").Value = [double]"-18.762619497260452"
$ws.Range("B4").Value = [double]"-5.1938099456774731E-2"
$ws.Range("A5").Value = [double]"-19.213902499303479"
$ws.Range("B5").Value = [double]"1.6215566659970464"
$ws.Range("A6").Value = [double]"-11.477784012965978"
$ws.Range("B6").Value = [double]"5.3253362240968904"
$ws.Range("A7").Value = [double]"-16.712062265707708"
$ws.Range("B7").Value = [double]"6.7620619070295112"
$ws.Range("A8").Value = [double]"17.398772667603374"
$ws.Range("B8").Value = [double]"12.665775030016647"
$ws.Range("A9").Value = [double]"-8.1211474688400553"
$ws.Range("B9").Value = [double]"1.721637871511049"
$ws.Range("C9").Value = [double]"-1"
$ws.Range("A10").Value = [double]"-14.019244026089597"
$ws.Range("B10").Value = [double]"-2.2027669800954284"
$ws.Range("A11").Value = [double]"-10.126233970702842"
$ws.Range("B11").Value = [double]"9.0914208275220822"
$ws.Range("D11").Value = [double]"1"
$ws.Range("A12").Value = [double]"7.1366443546039591"
$ws.Range("B12").Value = [double]"4.7560909738728592"
$ws.Range("A13").Value = [double]"-21.529648906424708"
$ws.Range("B13").Value = [double]"2.7339900409015634E-2"
$ws.Range("A14").Value = [double]"-24.861478598620529"
$ws.Range("B14").Value = [double]"-0.27924864813734385"
$ws.Range("A15").Value = [double]"-13.731266704883502"
$ws.Range("B15").Value = [double]"-3.5918967031508751"
$ws.Range("A16").Value = [double]"-8.086671755685007"
$ws.Range("B16").Value = [double]"-0.92779311623438798"
$ws.Range("C16").Value = [double]"0"
$ws.Range("A17").Value = [double]"-7.688470137958757"
$ws.Range("B17").Value = [double]"-1.8418377904994645"
$ws.Range("A18").Value = [double]"-24.327347629743567"
$ws.Range("B18").Value = [double]"-6.7529538521736017"
$ws.Range("A19").Value = [double]"-3.8104606250983761"
$ws.Range("B19").Value = [double]"12.581720944013441"
$ws.Range("A20").Value = [double]"6.0840805446344719"
$ws.Range("B20").Value = [double]"8.9854813352036764"
$ws.Range("C20").Value = [double]"1"
$ws.Range("A21").Value = [double]"13.252201117910749"
$ws.Range("B21").Value = [double]"17.344856514968459"
$ws.Range("A22").Value = [double]"-3.372349111463377"
$ws.Range("B22").Value = [double]"2.625542137375044"
$ws.Range("A23").Value = [double]"4.8480038015272848"
$ws.Range("B23").Value = [double]"10.597007347183315"
$ws.Range("A24").Value = [double]"-6.4292273870743779"
$ws.Range("B24").Value = [double]"3.2047844443032618"
$ws.Range("A25").Value = [double]"-11.711013854029947"
$ws.Range("B25").Value = [double]"-5.4303223761000243"
$ws.Range("A26").Value = [double]"-5.3934835600455315"
$ws.Range("B26").Value = [double]"5.0086905956888206"
$ws.Range("A27").Value = [double]"-7.1560825487814883"
$ws.Range("B27").Value = [double]"3.4796470349906397"
$ws.Range("A28").Value = [double]"-4.4291851525713142"
$ws.Range("B28").Value = [double]"7.2829315716837817"
$ws.Range("A29").Value = [double]"-2.4475031420629643"
$ws.Range("B29").Value = [double]"8.0613314205745734"
$ws.Range("A30").Value = [double]"-10.873827969576475"
$ws.Range("A31").Value = [double]"0.15864049147768533"
$ws.Range("B31").Value = [double]"11.764910001832355"
$ws.Range("C31").Value = [double]"1"
$ws.Range("A32").Value = [double]"-9.2299935727186195"
$ws.Range("B32").Value = [double]"3.3264225571310422"
$ws.Range("A33").Value = [double]"2.6501186201025959"
$ws.Range("B33").Value = [double]"12.794730333956959"
$ws.Range("A34").Value = [double]"-0.97112141811364827"
$ws.Range("B34").Value = [double]"7.2734026462298997"
$ws.Range("A35").Value = [double]"12.131513351561036"
$ws.Range("B35").Value = [double]"6.560894496417772"
$ws.Range("B36").Value = [double]"12.776229815279043"
$ws.Range("A37").Value = [double]"24.109819821928571"
$ws.Range("B37").Value = [double]"21.04196456646984"
$ws.Range("A38").Value = [double]"12.892348505638003"
$ws.Range("B38").Value = [double]"18.058104480505275"
$ws.Range("A39").Value = [double]"15.212012006084837"
$ws.Range("B39").Value = [double]"19.642951521202328"
$ws.Range("A40").Value = [double]"6.3153491417835195"
$ws.Range("B40").Value = [double]"7.4714656539394202"
$ws.Range("A41").Value = [double]"12.692016854258048"
$ws.Range("B41").Value = [double]"15.006505167144956"
$ws.Range("A42").Value = [double]"14.044924372030614"
$ws.Range("B42").Value = [double]"19.075535212858135"
$ws.Range("A43").Value = [double]"20.656866475143282"
$ws.Range("B43").Value = [double]"24.896528581236694"
$ws.Range("A44").Value = [double]"13.222956540584669"
$ws.Range("B44").Value = [double]"19.239993606155767"
$ws.Range("B45").Value = [double]"22.365475125576808"
$ws.Range("A46").Value = [double]"25.722624399712306"
$ws.Range("B46").Value = [double]"31.655256810818724"
$ws.Range("A47").Value = [double]"15.62119021778031"
$ws.Range("B47").Value = [double]"17.562304141326678"
$ws.Range("A48").Value = [double]"20.290157242676521"
$ws.Range("B48").Value = [double]"17.353504072624254"
$ws.Range("A49").Value = [double]"16.077147891735581"
$ws.Range("A50").Value = [double]"22.799589176657165"
$ws.Range("B50").Value = [double]"22.599613662957012"
$ws.Range("A51").Value = [double]"6.4082551977968301"
$ws.Range("B51").Value = [double]"9.3294471828599423"
$ws.Range("A52").Value = [double]"-24.873135113216243"
$ws.Range("B52").Value = [double]"-25.755349079564191"
$ws.Range("A53").Value = [double]"-24.875083265084267"
$ws.Range("B53").Value = [double]"-24.409758526275617"
$ws.Range("A54").Value = [double]"-37.544353187634499"
$ws.Range("B54").Value = [double]"-40.684975259707841"
$ws.Range("A55").Value = [double]"-10.778018334986911"
$ws.Range("B55").Value = [double]"-11.454154094149816"
$ws.Range("A56").Value = [double]"-11.136454681680419"
$ws.Range("B56").Value = [double]"-10.262091694722585"
$ws.Range("A57").Value = [double]"-17.765383152890397"
$ws.Range("B57").Value = [double]"-17.720475352420188"
$ws.Range("A58").Value = [double]"-30.792502924531966"
$ws.Range("B58").Value = [double]"-36.145824948859996"
$ws.Range("A59").Value = [double]"-12.511234904713865"
$ws.Range("B59").Value = [double]"-14.946722086268622"
$ws.Range("A60").Value = [double]"-8.9853182332310997"
$ws.Range("B60").Value = [double]"-8.8869956869158813"
$ws.Range("A61").Value = [double]"-21.398317412547868"
$ws.Range("B61").Value = [double]"-25.118466490239317"
$ws.Range("C61").Value = [double]"0"
$ws.Range("A62").Value = [double]"-24.873746951706174"
$ws.Range("A63").Value = [double]"-24.902638210545732"
$ws.Range("B63").Value = [double]"-26.397074514756589"
$ws.Range("A64").Value = [double]"-18.86208241199752"
$ws.Range("B64").Value = [double]"-18.881073905738333"
$ws.Range("A65").Value = [double]"-13.509683055749477"
$ws.Range("B65").Value = [double]"-13.939138596374963"
$ws.Range("A66").Value = [double]"-23.706394564714799"
$ws.Range("B66").Value = [double]"-25.736659939356962"
$ws.Range("A67").Value = [double]"-33.118021758364293"
$ws.Range("B67").Value = [double]"-27.973104672777847"
$ws.Range("A68").Value = [double]"-34.99258303305264"
$ws.Range("B68").Value = [double]"-27.42414430317514"
$ws.Range("A69").Value = [double]"-33.963113678167957"
$ws.Range("B69").Value = [double]"-25.050919946027729"
$ws.Range("A70").Value = [double]"-22.904701527763415"
$ws.Range("B70").Value = [double]"-17.921402496864111"
$ws.Range("A71").Value = [double]"-37.848725049947952"
$ws.Range("B71").Value = [double]"-35.613090562840448"
$ws.Range("A72").Value = [double]"-43.154281196093386"
$ws.Range("B72").Value = [double]"-42.711110919854995"
$ws.Range("A73").Value = [double]"-40.424150088831986"
$ws.Range("B73").Value = [double]"-35.495343115765799"
$ws.Range("A74").Value = [double]"-26.50770118653357"
$ws.Range("B74").Value = [double]"-17.706431315247762"
$ws.Range("A75").Value = [double]"-25.678959939580214"
$ws.Range("B75").Value = [double]"-22.332681922100271"
$ws.Range("A76").Value = [double]"-24.701195170738917"
$ws.Range("B76").Value = [double]"-17.750409923852693"
$ws.Range("A77").Value = [double]"-20.919435762154773"
$ws.Range("B77").Value = [double]"-14.709197735329703"
$ws.Range("A78").Value = [double]"-19.599885212518146"
$ws.Range("B78").Value = [double]"-16.485362718893704"
$ws.Range("A79").Value = [double]"-34.086871540602985"
$ws.Range("B79").Value = [double]"-23.866553450233123"
$ws.Range("A80").Value = [double]"-29.038292927411526"
$ws.Range("B80").Value = [double]"-18.864035407723804"
$ws.Range("A81").Value = [double]"-18.060895945245818"
$ws.Range("B81").Value = [double]"-5.0938708355221838"
$ws.Range("A82").Value = [double]"-24.504368395250665"
$ws.Range("B82").Value = [double]"-15.09641347107663"

$ws.Activate()
$ws.Range("A67:D82").Select()
